# Updates the cryptos price/volume table (and the Aave/BabyDogeCoin row-45/46 swap).
# Generated from the authoritative cell-level diff of the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell, new text value, and whether the value must be
# forced to Text format first (otherwise Excel would auto-parse a numeric-
# looking string such as "1.001" or "244.71" into a floating point number).
$updates = @(
    @{ Cell = "D2"; Value = '29.087.00'; ForceText = $false }
    @{ Cell = "E2"; Value = '  +0.19%  '; ForceText = $false }
    @{ Cell = "D3"; Value = '1.836.35'; ForceText = $false }
    @{ Cell = "E3"; Value = '  +0.43%  '; ForceText = $false }
    @{ Cell = "D4"; Value = '1.001'; ForceText = $true }
    @{ Cell = "E4"; Value = '  +0.16%  '; ForceText = $false }
    @{ Cell = "D5"; Value = '244.71'; ForceText = $true }
    @{ Cell = "E5"; Value = '  +1.66%  '; ForceText = $false }
    @{ Cell = "D6"; Value = '0.6361'; ForceText = $true }
    @{ Cell = "E6"; Value = '  +2.42%  '; ForceText = $false }
    @{ Cell = "D7"; Value = '1.002'; ForceText = $true }
    @{ Cell = "E7"; Value = '  +0.14%  '; ForceText = $false }
    @{ Cell = "D8"; Value = '0.07579'; ForceText = $true }
    @{ Cell = "E8"; Value = '  +1.57%  '; ForceText = $false }
    @{ Cell = "D9"; Value = '0.2951'; ForceText = $true }
    @{ Cell = "E9"; Value = '  +1.20%  '; ForceText = $false }
    @{ Cell = "D10"; Value = '22.94'; ForceText = $true }
    @{ Cell = "E10"; Value = '  +1.29%  '; ForceText = $false }
    @{ Cell = "D11"; Value = '0.07756'; ForceText = $true }
    @{ Cell = "E11"; Value = '  +1.98%  '; ForceText = $false }
    @{ Cell = "D12"; Value = '1.838.45'; ForceText = $false }
    @{ Cell = "E12"; Value = '  +0.52%  '; ForceText = $false }
    @{ Cell = "D13"; Value = '5.017'; ForceText = $true }
    @{ Cell = "E13"; Value = '  +1.28%  '; ForceText = $false }
    @{ Cell = "D14"; Value = '0.6724'; ForceText = $true }
    @{ Cell = "E14"; Value = '  +1.56%  '; ForceText = $false }
    @{ Cell = "D15"; Value = '83.33'; ForceText = $true }
    @{ Cell = "E15"; Value = '  +1.65%  '; ForceText = $false }
    @{ Cell = "D16"; Value = '0.000009589'; ForceText = $true }
    @{ Cell = "E16"; Value = '  +4.54%  '; ForceText = $false }
    @{ Cell = "D17"; Value = '6.127'; ForceText = $true }
    @{ Cell = "E17"; Value = '  +2.14%  '; ForceText = $false }
    @{ Cell = "D18"; Value = '29.112.20'; ForceText = $false }
    @{ Cell = "E18"; Value = '  +0.27%  '; ForceText = $false }
    @{ Cell = "D19"; Value = '12.62'; ForceText = $true }
    @{ Cell = "E19"; Value = '  +2.32%  '; ForceText = $false }
    @{ Cell = "D20"; Value = '227.33'; ForceText = $true }
    @{ Cell = "E20"; Value = '  +1.06%  '; ForceText = $false }
    @{ Cell = "D21"; Value = '1.001'; ForceText = $true }
    @{ Cell = "E21"; Value = '  +0.07%  '; ForceText = $false }
    @{ Cell = "D22"; Value = '7.233'; ForceText = $true }
    @{ Cell = "E22"; Value = '  +0.95%  '; ForceText = $false }
    @{ Cell = "D23"; Value = '1.002'; ForceText = $true }
    @{ Cell = "E23"; Value = '  +0.12%  '; ForceText = $false }
    @{ Cell = "D24"; Value = '160.88'; ForceText = $true }
    @{ Cell = "E24"; Value = '  +0.94%  '; ForceText = $false }
    @{ Cell = "D25"; Value = '0.1405'; ForceText = $true }
    @{ Cell = "E25"; Value = '  +3.85%  '; ForceText = $false }
    @{ Cell = "D26"; Value = '8.560'; ForceText = $true }
    @{ Cell = "E26"; Value = '  +1.95%  '; ForceText = $false }
    @{ Cell = "D27"; Value = '18.00'; ForceText = $true }
    @{ Cell = "E27"; Value = '  +1.21%  '; ForceText = $false }
    @{ Cell = "D28"; Value = '1.501'; ForceText = $true }
    @{ Cell = "E28"; Value = '  +0.36%  '; ForceText = $false }
    @{ Cell = "D29"; Value = '4.132'; ForceText = $true }
    @{ Cell = "E29"; Value = '  +2.05%  '; ForceText = $false }
    @{ Cell = "D30"; Value = '4.091'; ForceText = $true }
    @{ Cell = "E30"; Value = '  +1.63%  '; ForceText = $false }
    @{ Cell = "D31"; Value = '1.205'; ForceText = $true }
    @{ Cell = "E31"; Value = '  +0.27%  '; ForceText = $false }
    @{ Cell = "D32"; Value = '0.05421'; ForceText = $true }
    @{ Cell = "E32"; Value = '  +3.67%  '; ForceText = $false }
    @{ Cell = "D33"; Value = '1.867'; ForceText = $true }
    @{ Cell = "E33"; Value = '  +1.77%  '; ForceText = $false }
    @{ Cell = "D34"; Value = '0.7484'; ForceText = $true }
    @{ Cell = "E34"; Value = '  +1.91%  '; ForceText = $false }
    @{ Cell = "D35"; Value = '1.144'; ForceText = $true }
    @{ Cell = "E35"; Value = '  -0.34%  '; ForceText = $false }
    @{ Cell = "D36"; Value = '2.663'; ForceText = $true }
    @{ Cell = "E36"; Value = '  +0.63%  '; ForceText = $false }
    @{ Cell = "D37"; Value = '1.239.89'; ForceText = $false }
    @{ Cell = "E37"; Value = '  -2.84%  '; ForceText = $false }
    @{ Cell = "D39"; Value = '0.01792'; ForceText = $true }
    @{ Cell = "E39"; Value = '  +0.70%  '; ForceText = $false }
    @{ Cell = "D40"; Value = '6.630'; ForceText = $true }
    @{ Cell = "E40"; Value = '  +5.03%  '; ForceText = $false }
    @{ Cell = "D41"; Value = '0.9068'; ForceText = $true }
    @{ Cell = "E41"; Value = '  +1.68%  '; ForceText = $false }
    @{ Cell = "D42"; Value = '1.002'; ForceText = $true }
    @{ Cell = "E42"; Value = '  +0.17%  '; ForceText = $false }
    @{ Cell = "D43"; Value = '102.39'; ForceText = $true }
    @{ Cell = "E43"; Value = '  +0.70%  '; ForceText = $false }
    @{ Cell = "D44"; Value = '1.996.40'; ForceText = $false }
    @{ Cell = "E44"; Value = '  +0.96%  '; ForceText = $false }
    @{ Cell = "B45"; Value = 'Aave'; ForceText = $false }
    @{ Cell = "C45"; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; ForceText = $false }
    @{ Cell = "D45"; Value = '65.17'; ForceText = $true }
    @{ Cell = "E45"; Value = '  +2.80%  '; ForceText = $false }
    @{ Cell = "B46"; Value = 'BabyDogeCoin'; ForceText = $false }
    @{ Cell = "C46"; Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; ForceText = $false }
    @{ Cell = "D46"; Value = '0.00000000123'; ForceText = $true }
    @{ Cell = "E46"; Value = '  +2.92%  '; ForceText = $false }
    @{ Cell = "D47"; Value = '0.5119'; ForceText = $true }
    @{ Cell = "E47"; Value = '  +0.10%  '; ForceText = $false }
    @{ Cell = "D48"; Value = '0.4105'; ForceText = $true }
    @{ Cell = "E48"; Value = '  +3.68%  '; ForceText = $false }
    @{ Cell = "D49"; Value = '9.147'; ForceText = $true }
    @{ Cell = "E49"; Value = '  +3.49%  '; ForceText = $false }
    @{ Cell = "D50"; Value = '6.783'; ForceText = $true }
    @{ Cell = "E50"; Value = '  +1.96%  '; ForceText = $false }
    @{ Cell = "D51"; Value = '0.05781'; ForceText = $true }
    @{ Cell = "E51"; Value = '  +0.53%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
